# Generate Report for Handoff
#
# The handoff-status sample data is regenerated: the old placeholder
# files (random-guid .png/.md names) are swapped for the new
# calleeMd1.md / calleeMd2.md / callerMd1.md / callerMd2.md fixture set,
# handoff timestamps are refreshed, and a 4th data row (callerMd2.md) is
# added to every sheet (Overview, zh-cn, de-de) together with its
# matching hyperlinks.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$hyperlinkColor = 15570276   # OLE (BGR) form of RGB 6495ED, matches the workbook's HyperLink font color

function Set-Text($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
}

function Set-DateText($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
    $ws.Range($addr).NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

function Style-AsLink($ws, $addr) {
    $ws.Range($addr).Font.Underline = 2
    $ws.Range($addr).Font.Color = $hyperlinkColor
}

function Add-Hlink($ws, $addr, $url, $label) {
    $ws.Hyperlinks.Add($ws.Range($addr), $url, "", "", $label) | Out-Null
}

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
Set-Text      $ws1 "A2" "calleeMd1.md"
Set-Text      $ws1 "B2" "Ready for handoff"
Set-Text      $ws1 "C2" "Ready for handoff"
Set-DateText  $ws1 "D2" "2016-03-21 15:01:40"

Set-Text      $ws1 "A3" "calleeMd2.md"
Set-Text      $ws1 "B3" "Ready for handoff"
Set-Text      $ws1 "C3" "Ready for handoff"
Set-DateText  $ws1 "D3" "2016-03-21 15:01:40"

Set-Text      $ws1 "A4" "callerMd1.md"
Set-Text      $ws1 "B4" "Ready for handoff"
Set-Text      $ws1 "C4" "Ready for handoff"
Set-DateText  $ws1 "D4" "2016-03-21 15:01:40"

Set-Text      $ws1 "A5" "callerMd2.md"
Set-Text      $ws1 "B5" "Ready for handoff"
Set-Text      $ws1 "C5" "Ready for handoff"
Set-DateText  $ws1 "D5" "2016-03-21 15:01:40"

Style-AsLink $ws1 "A5"

$ws1.Hyperlinks.Delete()
Add-Hlink $ws1 "A2" "https://github.com/OpenLocalizationTest/oltest/blob/cc5b6e0f7d1b6d61a7f8bd72b466585f233987c2/e2e/calleeMd1.md" "calleeMd1.md"
Add-Hlink $ws1 "A3" "https://github.com/OpenLocalizationTest/oltest/blob/cc5b6e0f7d1b6d61a7f8bd72b466585f233987c2/e2e/calleeMd2.md" "calleeMd2.md"
Add-Hlink $ws1 "A4" "https://github.com/OpenLocalizationTest/oltest/blob/cc5b6e0f7d1b6d61a7f8bd72b466585f233987c2/e2e/callerMd1.md" "callerMd1.md"
Add-Hlink $ws1 "A5" "https://github.com/OpenLocalizationTest/oltest/blob/cc5b6e0f7d1b6d61a7f8bd72b466585f233987c2/e2e/callerMd2.md" "callerMd2.md"

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
# Row 4 used to carry a "Dependency From" value in K4; the new sample set
# records that relationship in the "Reference Tokens" column (I4)
# instead, so the stale K4 cell needs to be cleared out entirely.
$ws2.Range("K4").ClearContents()

Set-Text      $ws2 "A2" "calleeMd1.md"
Set-Text      $ws2 "B2" ".md"
Set-Text      $ws2 "C2" "Ready for handoff"
Set-Text      $ws2 "D2" "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf"
Set-DateText  $ws2 "E2" "2016-03-21 15:01:36"
Set-DateText  $ws2 "H2" "0001-01-01 00:00:00"
Set-Text      $ws2 "J2" "Include"
Set-Text      $ws2 "K2" "e2e\callerMd2.md,`ne2e\callerMd1.md"

Set-Text      $ws2 "A3" "calleeMd2.md"
Set-Text      $ws2 "B3" ".md"
Set-Text      $ws2 "C3" "Ready for handoff"
Set-Text      $ws2 "D3" "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf"
Set-DateText  $ws2 "E3" "2016-03-21 15:01:36"
Set-DateText  $ws2 "H3" "0001-01-01 00:00:00"
Set-Text      $ws2 "J3" "Include"
Set-Text      $ws2 "K3" "e2e\callerMd1.md"

Set-Text      $ws2 "A4" "callerMd1.md"
Set-Text      $ws2 "B4" ".md"
Set-Text      $ws2 "C4" "Ready for handoff"
Set-Text      $ws2 "D4" "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf"
Set-DateText  $ws2 "E4" "2016-03-21 15:01:36"
Set-DateText  $ws2 "H4" "0001-01-01 00:00:00"
Set-Text      $ws2 "I4" "e2e\calleeMd1.md,`ne2e\calleeMd2.md"
Set-Text      $ws2 "J4" "Include"

Set-Text      $ws2 "A5" "callerMd2.md"
Set-Text      $ws2 "B5" ".md"
Set-Text      $ws2 "C5" "Ready for handoff"
Set-Text      $ws2 "D5" "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf"
Set-DateText  $ws2 "E5" "2016-03-21 15:01:36"
Set-DateText  $ws2 "H5" "0001-01-01 00:00:00"
Set-Text      $ws2 "I5" "e2e\calleeMd1.md"
Set-Text      $ws2 "J5" "Include"

Style-AsLink $ws2 "A5"
Style-AsLink $ws2 "D5"

$ws2.Hyperlinks.Delete()
Add-Hlink $ws2 "A2" "https://github.com/OpenLocalizationTest/oltest/blob/cc5b6e0f7d1b6d61a7f8bd72b466585f233987c2/e2e/calleeMd1.md" "calleeMd1.md"
Add-Hlink $ws2 "D2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a8939cd7ad7f2ced89fc02704f27db975618dcea/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf" "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf"
Add-Hlink $ws2 "A3" "https://github.com/OpenLocalizationTest/oltest/blob/cc5b6e0f7d1b6d61a7f8bd72b466585f233987c2/e2e/calleeMd2.md" "calleeMd2.md"
Add-Hlink $ws2 "D3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a8939cd7ad7f2ced89fc02704f27db975618dcea/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf" "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf"
Add-Hlink $ws2 "A4" "https://github.com/OpenLocalizationTest/oltest/blob/cc5b6e0f7d1b6d61a7f8bd72b466585f233987c2/e2e/callerMd1.md" "callerMd1.md"
Add-Hlink $ws2 "D4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a8939cd7ad7f2ced89fc02704f27db975618dcea/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf" "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf"
Add-Hlink $ws2 "A5" "https://github.com/OpenLocalizationTest/oltest/blob/cc5b6e0f7d1b6d61a7f8bd72b466585f233987c2/e2e/callerMd2.md" "callerMd2.md"
Add-Hlink $ws2 "D5" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a8939cd7ad7f2ced89fc02704f27db975618dcea/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf" "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf"

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3.Range("K4").ClearContents()

Set-Text      $ws3 "A2" "calleeMd1.md"
Set-Text      $ws3 "B2" ".md"
Set-Text      $ws3 "C2" "Ready for handoff"
Set-Text      $ws3 "D2" "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf"
Set-DateText  $ws3 "E2" "2016-03-21 15:01:40"
Set-DateText  $ws3 "H2" "0001-01-01 00:00:00"
Set-Text      $ws3 "J2" "Include"
Set-Text      $ws3 "K2" "e2e\callerMd2.md,`ne2e\callerMd1.md"

Set-Text      $ws3 "A3" "calleeMd2.md"
Set-Text      $ws3 "B3" ".md"
Set-Text      $ws3 "C3" "Ready for handoff"
Set-Text      $ws3 "D3" "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf"
Set-DateText  $ws3 "E3" "2016-03-21 15:01:40"
Set-DateText  $ws3 "H3" "0001-01-01 00:00:00"
Set-Text      $ws3 "J3" "Include"
Set-Text      $ws3 "K3" "e2e\callerMd1.md"

Set-Text      $ws3 "A4" "callerMd1.md"
Set-Text      $ws3 "B4" ".md"
Set-Text      $ws3 "C4" "Ready for handoff"
Set-Text      $ws3 "D4" "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf"
Set-DateText  $ws3 "E4" "2016-03-21 15:01:40"
Set-DateText  $ws3 "H4" "0001-01-01 00:00:00"
Set-Text      $ws3 "I4" "e2e\calleeMd1.md,`ne2e\calleeMd2.md"
Set-Text      $ws3 "J4" "Include"

Set-Text      $ws3 "A5" "callerMd2.md"
Set-Text      $ws3 "B5" ".md"
Set-Text      $ws3 "C5" "Ready for handoff"
Set-Text      $ws3 "D5" "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf"
Set-DateText  $ws3 "E5" "2016-03-21 15:01:40"
Set-DateText  $ws3 "H5" "0001-01-01 00:00:00"
Set-Text      $ws3 "I5" "e2e\calleeMd1.md"
Set-Text      $ws3 "J5" "Include"

Style-AsLink $ws3 "A5"
Style-AsLink $ws3 "D5"

$ws3.Hyperlinks.Delete()
Add-Hlink $ws3 "A2" "https://github.com/OpenLocalizationTest/oltest/blob/cc5b6e0f7d1b6d61a7f8bd72b466585f233987c2/e2e/calleeMd1.md" "calleeMd1.md"
Add-Hlink $ws3 "D2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/59c805587fcb407a6443b20ff11cfc4d13ca2b1f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf" "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf"
Add-Hlink $ws3 "A3" "https://github.com/OpenLocalizationTest/oltest/blob/cc5b6e0f7d1b6d61a7f8bd72b466585f233987c2/e2e/calleeMd2.md" "calleeMd2.md"
Add-Hlink $ws3 "D3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/59c805587fcb407a6443b20ff11cfc4d13ca2b1f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf" "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf"
Add-Hlink $ws3 "A4" "https://github.com/OpenLocalizationTest/oltest/blob/cc5b6e0f7d1b6d61a7f8bd72b466585f233987c2/e2e/callerMd1.md" "callerMd1.md"
Add-Hlink $ws3 "D4" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/59c805587fcb407a6443b20ff11cfc4d13ca2b1f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf" "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf"
Add-Hlink $ws3 "A5" "https://github.com/OpenLocalizationTest/oltest/blob/cc5b6e0f7d1b6d61a7f8bd72b466585f233987c2/e2e/callerMd2.md" "callerMd2.md"
Add-Hlink $ws3 "D5" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/59c805587fcb407a6443b20ff11cfc4d13ca2b1f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf" "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf"

Write-Output "Report regenerated for handoff"
